# "Added print graph correlations"
# - Update the Sheet1 symbol list (A1:A9) with the new set of tickers.
# - Add a new "Correlations" worksheet (placed after Sheet1) holding the
#   4x4 correlation matrix for those tickers.
# - Leave Sheet1 as the active sheet with A9 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Sheet1: refresh the symbols column -------------------------------
$ws.Range("A1").Value = "Symbols"
$ws.Range("A2").Value = "AMZN"
$ws.Range("A3").Value = "COST"
$ws.Range("A4").Value = "GOOD"
$ws.Range("A5").Value = "SQ"
$ws.Range("A6").Value = "TSLA"
$ws.Range("A7").Value = "MDB"
$ws.Range("A8").Value = "BYND"
$ws.Range("A9").Value = "TCS.TO"

# --- Add the Correlations sheet right after Sheet1 ---------------------
$corr = $wb.Worksheets.Add($null, $ws)
$corr.Name = "Correlations"

$corr.Range("A1").Value = 1
$corr.Range("B1").Value = -0.022882770093682309
$corr.Range("C1").Value = 0.097151317164047138
$corr.Range("D1").Value = 0.67154125428006239

$corr.Range("A2").Value = -0.022882770093682309
$corr.Range("B2").Value = 1
$corr.Range("C2").Value = 0.54359152388547205
$corr.Range("D2").Value = 0.37389212624620721

$corr.Range("A3").Value = 0.097151317164047124
$corr.Range("B3").Value = 0.54359152388547216
$corr.Range("C3").Value = 1
$corr.Range("D3").Value = 0.27242058256347218

$corr.Range("A4").Value = 0.67154125428006228
$corr.Range("B4").Value = 0.37389212624620721
$corr.Range("C4").Value = 0.27242058256347218
$corr.Range("D4").Value = 1

# --- Restore Sheet1 as the active sheet/selection -----------------------
$ws.Activate() | Out-Null
$ws.Range("A9").Select() | Out-Null
